# CORE_holdings.xlsx update
#   1. Bump the "as of" date in the confidential disclosure footnote (A11)
#      from 2021-04-09 to 2021-04-21.
#   2. Refresh the Weight (D) and Percent Change (E) figures for the model
#      holdings rows (2-8) with the latest snapshot values.
#
# The sheet ships protected, so it has to be unprotected before the edits
# and re-protected afterwards to restore the original state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- 1. Update the disclosure footnote's "as of" date ----------------------
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."
$ws.Range("A11").Value = $newText

# --- 2. Refresh Weight / Percent Change figures -----------------------------
$ws.Range("D2").Value = 0.4922177157500554
$ws.Range("E2").Value = 0.01147107897664079

$ws.Range("D3").Value = 0.252544514914781
$ws.Range("E3").Value = 0.007381676074685206

$ws.Range("D4").Value = 0.09836056011603624
$ws.Range("E4").Value = 0.01706102117061037

$ws.Range("D5").Value = 0.1008414961833248
$ws.Range("E5").Value = 0.01744914682348409

$ws.Range("D6").Value = 0.02890796038279522
$ws.Range("E6").Value = 0.02128732849071824

$ws.Range("D7").Value = 0.02712775265300726
$ws.Range("E7").Value = 0.02317880794701987

$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = 0.01219236198334883

# --- 3. Restore sheet protection --------------------------------------------
$ws.Protect("D382")
